# Debuging incorrect values on output
# Expand distance bus and add debug hooks.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- "distance bus" debug table: id (A), distance value (B), squared (C) ---
$ids = @(17, 13, 10, 11, 48, 26, 60, 56, 62, 16)
$distances = @(
    2378858632.4404702,
    2546340798.0274501,
    2671753018.757,
    2894561484.2773299,
    2918228302.5657401,
    3089562004.6875601,
    3090682997.1241899,
    3109507391.0303402,
    3123623195.41642,
    3125181049.6946502
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $r = $i + 1
    $ws2.Range("A$r").Value = $ids[$i]
    $ws2.Range("B$r").Value = $distances[$i]
}

# six-decimal debug format on the raw distance column
$ws2.Range("B1:B10").NumberFormat = "0.000000"

# debug hook: square the distance (first/last rows entered individually,
# the middle filled as one block so Excel keeps them as a shared formula)
$ws2.Range("C1").Formula = "=POWER(B1,2)"
$ws2.Range("C2:C9").Formula = "=POWER(B2,2)"
$ws2.Range("C10").Formula = "=POWER(B10,2)"

# widen the new columns so the debug values are readable
$ws2.Columns.Item(2).ColumnWidth = 19.8
$ws2.Columns.Item(3).ColumnWidth = 27.1

$ws2.PageSetup.Orientation = 1

# make Sheet2 the active tab with C10 selected
$ws2.Activate()
$ws2.Range("C10").Select()
